$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# The starting point of R0 (row 4) used to be a fixed 0; it now continues
# the incrementing series started in A3, and the series below (previously
# starting at A5) shifts down to follow from the new A4.
$ws.Range("A4").Formula = "=A3+1"
$ws.Range("A5").Formula = "=A4+1"

# I3 referenced A5 plus the old offset from A4 (1) plus 3; now that A4 is
# part of the series, the formula only needs the remaining +3 offset.
$ws.Range("I3").Formula = "=A5+3"

[void]$ws.Range("I3").Select()
